$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("combined")

$ws.Range("D28").Value = 10
$ws.Range("E28").Value = 11
$ws.Range("F28").Value = 12
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = 14
$ws.Range("I28").Value = 15

$ws.Range("C29").Value = 0.001
$ws.Range("D29").Value = 0.8288288288288288
$ws.Range("E29").Value = 0.81981981981981977
$ws.Range("F29").Value = 0.83783783783783783
$ws.Range("G29").Value = 0.8288288288288288
$ws.Range("H29").Value = 0.83783783783783783
$ws.Range("I29").Value = 0.84684684684684686

$ws.Range("C30").Value = 0.01
$ws.Range("D30").Value = 0.8288288288288288
$ws.Range("E30").Value = 0.8288288288288288
$ws.Range("F30").Value = 0.81981981981981977
$ws.Range("G30").Value = 0.8288288288288288
$ws.Range("H30").Value = 0.8288288288288288
$ws.Range("I30").Value = 0.80180180180180183

$ws.Range("C31").Value = 0.05
$ws.Range("D31").Value = 0.70270270270270274
$ws.Range("E31").Value = 0.67567567567567566
$ws.Range("F31").Value = 0.68468468468468469
$ws.Range("G31").Value = 0.67567567567567566
$ws.Range("H31").Value = 0.68468468468468469
$ws.Range("I31").Value = 0.66666666666666663

$ws.Range("D28:I28").Font.Bold = $true
$ws.Range("C29").Font.Bold = $true
$ws.Range("C30").Font.Bold = $true
$ws.Range("C31").Font.Bold = $true
